$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, shifting existing rows 122:224 down to 123:225
$ws.Rows("122:122").Insert()

# Populate the newly inserted row 122 with the new price record
$ws.Range("A122").Value = 7
$ws.Range("B122").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C122").Value = "Ñuble"
$ws.Range("D122").Value = 44512
$ws.Range("D122").NumberFormat = $ws.Range("D121").NumberFormat
$ws.Range("E122").Value = 16
$ws.Range("F122").Value = 100114001
$ws.Range("G122").Value = "Papa"
$ws.Range("H122").Value = "Patagonia"
$ws.Range("I122").Value = "1a (guarda)"
$ws.Range("J122").Value = 200
$ws.Range("K122").Value = 7000
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = 7250
$ws.Range("N122").Value = "$/saco 25 kilos"
$ws.Range("O122").Value = "Provincia de Diguillín"
$ws.Range("P122").Value = 290
$ws.Range("Q122").Value = 25
$ws.Range("R122").Value = "Hortaliza"
